$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47, shifting existing rows 47:61 down to 48:62.
$ws.Rows("47:47").Insert()

# Populate the newly inserted row 47 with the new weekly price record.
$ws.Range("A47").Value = 10
$ws.Range("B47").Value = "Vega Modelo de Temuco"
$ws.Range("C47").Value = "La Araucanía"
$ws.Range("D47").Value = 44524
$ws.Range("E47").Value = 9
$ws.Range("F47").Value = "Fruta"
$ws.Range("G47").Value = 100101
$ws.Range("H47").Value = "Berries"
$ws.Range("I47").Value = 100101001
$ws.Range("J47").Value = "Arándano (blue)"
$ws.Range("K47").Value = "Sin especificar"
$ws.Range("L47").Value = "Primera"
$ws.Range("M47").Value = 300
$ws.Range("N47").Value = 3000
$ws.Range("O47").Value = 3000
$ws.Range("P47").Value = 3000
$ws.Range("Q47").Value = "$/kilo"
$ws.Range("R47").Value = "Región del Maule"
$ws.Range("S47").Value = 3000
$ws.Range("T47").Value = 1
